$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.096050500869751
$ws.Range("B1").Value = 1.968271732330322
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.020567178726196
$ws.Range("E1").Value = 1.128741383552551
